# Update cryptocurrency price/volume snapshot values (Price = column D, Volume(1h) = column E)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.673.84"
$ws.Range("E2").Value = "  -0.78%  "
$ws.Range("D3").Value = "1.886.84"
$ws.Range("E3").Value = "  -0.92%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  +0.21%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "236.42"
$ws.Range("E5").Value = "  -4.11%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.003"
$ws.Range("E6").Value = "  +0.34%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4879"
$ws.Range("E7").Value = "  -2.53%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2909"
$ws.Range("E8").Value = "  -2.86%  "
$ws.Range("E9").Value = "  -2.82%  "
$ws.Range("D10").Value = "1.888.73"
$ws.Range("E10").Value = "  -0.82%  "
$ws.Range("E11").Value = "  -4.67%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07254"
$ws.Range("E12").Value = "  -1.18%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "89.17"
$ws.Range("E13").Value = "  -2.79%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.013"
$ws.Range("E14").Value = "  -2.15%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6600"
$ws.Range("E15").Value = "  -3.33%  "
$ws.Range("D16").Value = "30.603.58"
$ws.Range("E16").Value = "  -0.96%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000007855"
$ws.Range("E17").Value = "  -2.81%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.003"
$ws.Range("E18").Value = "  +0.21%  "
$ws.Range("E19").Value = "  -3.18%  "
$ws.Range("D20").Value = "2.134.39"
$ws.Range("E20").Value = "  -0.64%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.005"
$ws.Range("E21").Value = "  +0.37%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.733"
$ws.Range("E22").Value = "  -2.95%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "190.53"
$ws.Range("E23").Value = "  +3.35%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.091"
$ws.Range("E24").Value = "  -0.51%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.297"
$ws.Range("E25").Value = "  -0.95%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "159.33"
$ws.Range("E26").Value = "  +3.43%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.27"
$ws.Range("E27").Value = "  -2.28%  "
$ws.Range("E28").Value = "  -6.35%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.407"
$ws.Range("E29").Value = "  +0.69%  "
$ws.Range("E30").Value = "  -3.50%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.09001"
$ws.Range("E31").Value = "  +0.04%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.929"
$ws.Range("E32").Value = "  -3.59%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05159"
$ws.Range("E33").Value = "  -2.51%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7254"
$ws.Range("E34").Value = "  -3.34%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.080"
$ws.Range("E35").Value = "  -5.57%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.698"
$ws.Range("E36").Value = "  +0.25%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.01810"
$ws.Range("E37").Value = "  -6.04%  "
$ws.Range("E38").Value = "  -2.25%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.9205"
$ws.Range("E39").Value = "  -2.51%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.047"
$ws.Range("E40").Value = "  -6.90%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.4376"
$ws.Range("E41").Value = "  -0.79%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "104.70"
$ws.Range("E42").Value = "  -1.76%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.9989"
$ws.Range("E43").Value = "  -0.10%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.713"
$ws.Range("E44").Value = "  -2.38%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.1328"
$ws.Range("E45").Value = "  -2.46%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "7.331"
$ws.Range("E46").Value = "  -5.99%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.4044"
$ws.Range("E47").Value = "  +2.65%  "
$ws.Range("E48").Value = "  -0.35%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.670"
$ws.Range("E49").Value = "  +0.40%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.407"
$ws.Range("E50").Value = "  +0.78%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "33.23"
$ws.Range("E51").Value = "  -0.87%  "
